# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1076
    5  = 3074
    7  = 2356
    11 = 1155
    13 = 51
    15 = 890
    16 = 284
    17 = 304
    19 = 17
    20 = 100
    21 = 57
    22 = 76
    23 = 2
    24 = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
